$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update match rows 167-171 with refreshed league-base data (id, date, teams, odds).
# The "id" column (B) holds numeric-looking strings; preserve the original cell
# style/format (General, no quote-prefix) while forcing text storage, by saving
# and restoring .Style around the apostrophe-prefixed assignment.

# Row 167
$origStyle = $ws.Range("B167").Style
$ws.Range("B167").Value2 = "'7060686"
$ws.Range("B167").Style = $origStyle
$ws.Range("D167").Value2 = 45429.54166666666
$ws.Range("E167").Value2 = "Istra 1961"
$ws.Range("F167").Value2 = "NK Varazdin"
$ws.Range("J167").Value2 = 2.2
$ws.Range("K167").Value2 = 3.1
$ws.Range("L167").Value2 = 3.3
$ws.Range("M167").Value2 = 2.2
$ws.Range("N167").Value2 = 3.1
$ws.Range("O167").Value2 = 3.3
$ws.Range("P167").Value2 = -0.25
$ws.Range("Q167").Value2 = 1.925
$ws.Range("R167").Value2 = 1.925
$ws.Range("S167").Value2 = 2.25
$ws.Range("T167").Value2 = 1.975
$ws.Range("U167").Value2 = 1.875

# Row 168
$origStyle = $ws.Range("B168").Style
$ws.Range("B168").Value2 = "'7097427"
$ws.Range("B168").Style = $origStyle
$ws.Range("D168").Value2 = 45430.5
$ws.Range("E168").Value2 = "NK Rudes"
$ws.Range("F168").Value2 = "NK Lokomotiva Zagreb"
$ws.Range("J168").Value2 = 8
$ws.Range("K168").Value2 = 4.75
$ws.Range("L168").Value2 = 1.333
$ws.Range("M168").Value2 = 8
$ws.Range("N168").Value2 = 4.75
$ws.Range("O168").Value2 = 1.333
$ws.Range("P168").Value2 = 1.5
$ws.Range("Q168").Value2 = 1.825
$ws.Range("R168").Value2 = 2.025
$ws.Range("S168").Value2 = 3
$ws.Range("T168").Value2 = 2.025
$ws.Range("U168").Value2 = 1.825

# Row 169
$origStyle = $ws.Range("B169").Style
$ws.Range("B169").Value2 = "'7093590"
$ws.Range("B169").Style = $origStyle
$ws.Range("D169").Value2 = 45430.59027777778
$ws.Range("E169").Value2 = "Slaven Belupo"
$ws.Range("F169").Value2 = "Dinamo Zagreb"
$ws.Range("J169").Value2 = 6.5
$ws.Range("K169").Value2 = 5
$ws.Range("L169").Value2 = 1.363
$ws.Range("M169").Value2 = 6.5
$ws.Range("N169").Value2 = 5
$ws.Range("O169").Value2 = 1.363
$ws.Range("P169").Value2 = 1.25
$ws.Range("Q169").Value2 = 2
$ws.Range("R169").Value2 = 1.85
$ws.Range("S169").Value2 = 2.75
$ws.Range("T169").Value2 = 1.875
$ws.Range("U169").Value2 = 1.975

# Row 170
$origStyle = $ws.Range("B170").Style
$ws.Range("B170").Value2 = "'7097426"
$ws.Range("B170").Style = $origStyle
$ws.Range("D170").Value2 = 45431.5
$ws.Range("E170").Value2 = "Hajduk Split"
$ws.Range("F170").Value2 = "HNK Gorica"
$ws.Range("J170").Value2 = 1.285
$ws.Range("K170").Value2 = 5
$ws.Range("L170").Value2 = 9.5
$ws.Range("M170").Value2 = 1.285
$ws.Range("N170").Value2 = 5
$ws.Range("O170").Value2 = 9.5
$ws.Range("P170").Value2 = -1.5
$ws.Range("Q170").Value2 = 1.875
$ws.Range("R170").Value2 = 1.975
$ws.Range("S170").Value2 = 3
$ws.Range("T170").Value2 = 2
$ws.Range("U170").Value2 = 1.85

# Row 171
$origStyle = $ws.Range("B171").Style
$ws.Range("B171").Value2 = "'7098238"
$ws.Range("B171").Style = $origStyle
$ws.Range("D171").Value2 = 45431.60416666666
$ws.Range("E171").Value2 = "NK Osijek"
$ws.Range("F171").Value2 = "HNK Rijeka"
$ws.Range("J171").Value2 = 1.95
$ws.Range("K171").Value2 = 3.5
$ws.Range("L171").Value2 = 3.5
$ws.Range("M171").Value2 = 1.9
$ws.Range("N171").Value2 = 3.6
$ws.Range("O171").Value2 = 3.7
$ws.Range("P171").Value2 = -0.5
$ws.Range("Q171").Value2 = 1.9
$ws.Range("R171").Value2 = 1.95
$ws.Range("S171").Value2 = 2.75
$ws.Range("T171").Value2 = 2.025
$ws.Range("U171").Value2 = 1.825
